$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty row 8 with new risk entry data
$ws.Range("A8").Value = "Week2.4"
$ws.Range("B8").Value = "Duplicated cards"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "use a loop and iterate only once to prevent it from happening"

# Move the active selection from F7 to C10 (reflects where the user clicked next)
[void]$ws.Range("C10").Select()
